$d = $word.ActiveDocument

# Update the portfolio URL text (hyperlink display text) from the old
# heroku app address to the new mohmedvaid.com address. The hyperlink
# target itself (rId8) is left untouched by the source diff.
$d.Content.Find.Execute(
    "https://mohmedvaid.herokuapp.com/", $false, $false, $false, $false,
    $false, $true, 1, $false, "https://mohmedvaid.com", 2
)
